# Update cryptocurrency price/volume figures per the latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.899.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.516.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.69%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.73"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.38"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.37%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.516.32"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.68%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.90%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.70%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.969.69"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.650.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.79%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.519.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.80"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.02%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.644.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0909"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "467.23"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.44%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.88%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.120"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.11"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.99"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.46"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.71%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.61"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.12%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -14.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.88%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.00%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.55%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0736"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.39%  "
